$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Clinical considerations" text for the first two turtles
# (age-category correction in the metadata):
#   Row 2 (TurtleID ID010 / MERRY FISHER): "Adult female" -> "Subadult female"
#   Row 3 (TurtleID ID047 / ZAL):          "Subadult ND"  -> "Juvenile ND"
$old2 = $ws.Range("C2").Value2
$old3 = $ws.Range("C3").Value2

$new2 = $old2 -replace "2019 HR Adult female", "2019 HR Subadult female"
$new3 = $old3 -replace "2019 HR Subadult ND", "2019 HR Juvenile ND"

# Write C3 first, then C2, so that newly appended shared-string entries
# end up in the same order as the reference workbook.
$ws.Range("C3").Value = $new3
$ws.Range("C2").Value = $new2

# Row 2's text grew longer, so its wrapped row height increases.
$ws.Rows.Item(2).RowHeight = 85

# Move the active selection to C2.
$ws.Range("C2").Select()
